$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: request_id -> task_id
$ws.Range("F1").Value = "task_id"

# Row 2 data updates (fal.ai -> kie.ai video generation record)
$ws.Range("A2").Value = "Meteorologist woman narrowly escapes tornado while reporting live"
$ws.Range("B2").Value = "Tornado chase live: meteorologist’s heart-stopping twist! 🌪️ #weatherwatch #stormchasing #tornadoalert #meteorology #viral #trending #livebroadcast #breakingnews #stormseason #extremeweather #newsflash #onair"
$ws.Range("C2").Value = "Open plains with dark swirling clouds, tornado funnel forming, meteorologist running with camera crew, cinematic realism, dynamic wind effects, intense sky contrast, emergency vehicles in background, handheld docu-style shot, suspenseful atmosphere, dust and debris flying around, lightning flashes intermittently in the stormy sky, dramatic close-ups on the meteorologist’s determined face, vivid weather phenomena dominating the scene, wide shots capturing the scale of the tornado and the chaos it causes, real-time live broadcast graphics overlayed on screen, urgency and tension palpable in the air, natural sounds of wind howling mixed with urgent voiceover, gritty realism with high-definition clarity, emotionally charged and immersive visual experience, blending raw nature power and human bravery, cinematic color grading emphasizing the storm’s ominous hues, fast-paced editing capturing every moment of the chase, immersive audio capturing the roar of the tornado and meteorologist’s breathless commentary.  "
$ws.Range("D2").Value = "A woman with wind-tossed hair and a soaked jacket holds a handheld camera close, her face streaked with rain and dirt as a massive tornado funnel spirals ominously behind her across the flat, open plains. She gasps into the mic, voice tight with adrenaline, “This thing’s coming faster than we thought—everyone, get clear, now!” The frame jitters as she stumbles forward, the camera shaking from her hurried steps and the howling wind that whips dust and debris past the lens. Lightning flashes illuminate the angry sky, casting sharp shadows on her determined expression as emergency vehicles flash red and blue lights in the distance. The shot is tight, selfie-style with a wide-angle lens that captures the swirling clouds pressing down and the tornado’s terrifying scale, filmed on high-definition stock with gritty, cinematic color grading emphasizing dark stormy blues and ominous grays. The relentless roar of the storm mingles with her breathless voiceover, the chaotic energy palpable as she glances nervously over her shoulder while pressing forward against the gale. Time of Day: late afternoon dusk. Lens: wide-angle handheld, high-definition stock. Audio: howling wind, crackling thunder, distant sirens. Background: swirling storm clouds, tornado funnel, emergency vehicles racing across the plains behind her."
$ws.Range("E2").Value = "in_progress"
$ws.Range("F2").Value = "e4e5cbd3f1c65f6aeeaa8b4bdb74ba00"
$ws.Range("G2").Value = "https://tempfile.aiquickdraw.com/p/a78a5a1538a444519c4b12d87c9c2bed_1752142446.mp4"
$ws.Range("I2").Value = "2025-07-10 12:10"

# Header row font color: theme color -> explicit black RGB
$ws.Range("A1:I1").Font.Color = 0

